# Updated convert_excel_date to throw ValueError when input- column values are
# of type non-numeric. To exercise this in the workbook, a new "Hire Date Str"
# column (L) is added that mirrors the "Hire Date" (E) column's serial-number
# values, but with a couple of rows replaced by non-numeric "dirty" values
# (subject names typed in by mistake) so the cleaning function has something
# to reject.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header: L1 = "Hire Date Str" -------------------------------------
# Give it the same (bold + bottom border) header formatting as the other
# header cells by copying K1's format over before writing the text.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = "Hire Date Str"

# --- Body cells: L2:L14 (row 9 is a deliberately blank spacer row, like E9) -
# Reuse the "Hire Date" (E) column's number format for the new column.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("L2:L8").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("L10:L14").PasteSpecial(-4122) | Out-Null

$ws.Range("L2").Value = 39690
$ws.Range("L3").Value = 39690
$ws.Range("L4").Value = 37118
$ws.Range("L6").Value = 41431
# Dirty, non-numeric values (order matters so new shared strings land in the
# same order as the source workbook: "physics" before "maths").
$ws.Range("L7").Value = "physics"
$ws.Range("L5").Value = "maths"
$ws.Range("L8").Value = 11037
$ws.Range("L10").Value = 32994
$ws.Range("L11").Value = 27919
$ws.Range("L12").Value = 42221
$ws.Range("L13").Value = 34700
$ws.Range("L14").Value = 40071

# Mirror the author's final selection/cursor position.
$ws.Range("L2").Select() | Out-Null
$excel.CutCopyMode = $false
